$d = $word.ActiveDocument

# --- 1. SSS-00003: insert "deve " before "permitir" ---
$d.Content.Find.Execute(
    "ao habilitar agendamento de visita, permitir que o atendente",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ao habilitar agendamento de visita, deve permitir que o atendente",
    2) | Out-Null

# --- 2. SSS-00004: shorten text, then remove the 3 bullet paragraphs that followed it ---
$d.Content.Find.Execute(
    "SSS-00004: O sistema deve conter um formulário para cadastro do cliente que solicitou o serviço, contendo as seguintes informações:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SSS-00004: O sistema deve conter um formulário para cadastro do cliente que solicitou o serviço.",
    2) | Out-Null

# Delete the bullet-list paragraphs ("Nome;", "Endereço da solicitação;", "Telefone para contato.")
# They always directly follow the SSS-00004 paragraph, so find it again and delete the next 3 paragraphs.
$p = $d.Paragraphs.Item(1)
while ($p.Range.Text -notmatch "^SSS-00004") {
    $p = $p.Next()
}
$bulletsStart = $p.Next().Range.Start
$bulletsEnd = $p.Next().Next().Next().Range.End
$d.Range($bulletsStart, $bulletsEnd).Delete() | Out-Null

# --- 3. SSS-00005 ---
$d.Content.Find.Execute(
    "SSS-00005: O sistema deve listar as matérias primas disponíveis em estoque.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SSS-00005: O sistema deve mostrar as matérias primas disponíveis em estoque filtrando por quantidade ou nome.",
    2) | Out-Null

# --- 4. SSS-00008 ---
$d.Content.Find.Execute(
    "SSS-00008: O sistema deve conter um campo de desconto opcional e um campo de taxa de locomoção obrigatórios no formulário de orçamento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SSS-00008: O sistema deve fazer a gestão das matérias primas do estoque da empresa.",
    2) | Out-Null

# --- 5. SSS-000015 ---
$d.Content.Find.Execute(
    "SSS-000015: O sistema deve disponibilizar todos as solicitações de serviço.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SSS-000015: O sistema deve disponibilizar todas as solicitações de serviço com opções de filtragem por status.",
    2) | Out-Null

# --- 6. SSS-000018 ---
$d.Content.Find.Execute(
    "SSS-000018: O sistema deve permitir que, caso o usuário escolha cancelar alguma solicitação, registre o motivo do cancelamento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SSS-000018: O sistema deve permitir que, caso o usuário escolha cancelar alguma solicitação, registre o motivo do cancelamento no campo utilizado para detalhes.",
    2) | Out-Null

# --- 7. Tail: replace the trailing empty paragraphs after SSS-000018 with new requirement paragraphs ---
# Locate the SSS-000018 paragraph again.
$p18 = $d.Paragraphs.Item(1)
while ($p18.Range.Text -notmatch "^SSS-000018") {
    $p18 = $p18.Next()
}
# keep one blank paragraph right after it, then remove everything else up to (and including) the
# last paragraph of the document.
$keepBlank = $p18.Next()
$tailStart = $keepBlank.Next().Range.Start
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailEnd = $lastPara.Range.End
$d.Range($tailStart, $tailEnd).Delete() | Out-Null

# Now append the new paragraphs after the kept blank paragraph.
$anchor = $keepBlank.Range
$anchor.InsertParagraphAfter()
$anchor.Collapse(0) | Out-Null

function Add-Paragraph([string]$text) {
    $script:anchor.InsertParagraphAfter()
    $p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $p.Range.Text = $text
    $script:anchor = $d.Paragraphs.Item($d.Paragraphs.Count).Range
}

$lastP = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastP.Range.Text = "SSS-000019: O website de marketing da empresa deve retratar identidade da empresa, a retratando visualmente e evidenciando os serviços que presta."
$lastP.Range.InsertParagraphAfter()
$lastP2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$lastP2.Range.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "SSS-000020: O sistema deve permitir que o usuário altere a quantidade de produtos de  determinada matéria prima no estoque."
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p2.Range.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "SSSS-000021: O sistema deve possibilitar o acesso ao sistema somente para usuários previamente cadastrados."
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$p2.Range.InsertParagraphAfter()

$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "SSS-000022: O sistema deve fazer a gestão das solicitações de serviço durante todo o ciclo de vida do pedido."
